$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.25
$ws.Range("I2").Value = 3.9
$ws.Range("J2").Value = 3.1
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("X2").Value = 9
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 26
$ws.Range("AC2").Value = 5
$ws.Range("AF2").Value = 101
$ws.Range("AG2").Value = 7.5
$ws.Range("AR2").Value = 101
$ws.Range("AZ2").Value = 81

# Row 3
$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 2.38
$ws.Range("Q3").Value = 2.7
$ws.Range("R3").Value = 1.44

# Row 5
$ws.Range("Q5").Value = 1.89
$ws.Range("R5").Value = 1.84

# Row 6
$ws.Range("Q6").Value = 1.84
$ws.Range("R6").Value = 1.89

# Row 10
$ws.Range("K10").Value = 1.91

# Row 11
$ws.Range("G11").Value = 2.15
$ws.Range("H11").Value = 3.25
$ws.Range("I11").Value = 3.5
$ws.Range("Q11").Value = 1.95
$ws.Range("R11").Value = 1.9

# Row 12
$ws.Range("J12").Value = 1.73
$ws.Range("Q12").Value = 1.93
$ws.Range("R12").Value = 1.93
